$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New daily auction records for 14-10-2021 and 15-10-2021
$ws.Range("A71").Value = "14-10-2021"
$ws.Range("B71").Value = 100000
$ws.Range("C71").Value = 140000
$ws.Range("D71").Value = 50000
$ws.Range("E71").Value = 40000
$ws.Range("F71").Value = 10000
$ws.Range("G71").Value = 2.85

$ws.Range("A72").Value = "15-10-2021"
$ws.Range("B72").Value = 100000
$ws.Range("C72").Value = 202000
$ws.Range("D72").Value = 100000
$ws.Range("E72").Value = 100000
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 2.84
